$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text-valued cells stay text (avoid Excel auto-converting
# numeric-looking strings like "5.200" or "0.5050" to numbers and
# dropping significant trailing zeros).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.624.21"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.01%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.841.25"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.16%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.38"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.19%  "

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.07%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4259"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.60%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.73%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07294"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.61%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8754"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.71%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.59"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.12%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.817.19"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.61%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.323"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.45%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.488"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.34%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06979"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.62%  "

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.12%  "

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.63%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008928"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.90%  "

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.13%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.32"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.87%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.637.95"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.11%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.959"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.38%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.32"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -2.33%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.064.52"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.39%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.995"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.49%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "155.57"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.69%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.45"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.52%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "119.34"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.09%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.200"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.62%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.866"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.40%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08857"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.50%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7566"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.77%  "

$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.494"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.46%  "

$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.937"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.35%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.127"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +2.42%  "

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.07%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05415"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.42%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.103"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.62%  "

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.05%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.812"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.10%  "

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.62%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5050"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.23%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.532"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -4.94%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.367"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.30%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.06547"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.91%  "

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.07%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "105.83"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.10%  "

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.07%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4621"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.69%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.632"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.23%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "64.21"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.46%  "
